$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears in column C, rows 2 and 3, on both the zh-cn and de-de sheets)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsZh = $wb.Worksheets("zh-cn")
$wsDe = $wb.Worksheets("de-de")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) -- handback has now happened
# ---------------------------------------------------------------------
$wsZh.Range("H2").Value = "2016-03-18 08:12:20"
$wsZh.Range("H3").Value = "2016-03-18 08:12:20"

$wsDe.Range("H2").Value = "2016-03-18 08:12:25"
$wsDe.Range("H3").Value = "2016-03-18 08:12:25"

# ---------------------------------------------------------------------
# 3. New columns populated as part of the handback: F = "Latest Target
#    File", G = "Latest Handback File". Both are hyperlinked file names,
#    matching the style already used for the other link columns (blue,
#    underlined).
# ---------------------------------------------------------------------
$hlColor = 15570276   # BGR long for RGB FF6495ED, matches existing HyperLink style

function Add-HandbackLink($ws, $cellRef, $address, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText)
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hlColor
}

# --- zh-cn sheet ---
Add-HandbackLink $wsZh "F2" "https://github.com/OpenLocalizationTest/oltest/blob/ab42d9a503830687685b01245f3dbd9aa1fb2bbb/e2e/0dab1796-b58b-47e1-aafe-a447b6b360bb.md" "0dab1796-b58b-47e1-aafe-a447b6b360bb.md"
Add-HandbackLink $wsZh "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4415d4333c52fc7d83a3c74ab476dfc18f340275/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/0dab1796-b58b-47e1-aafe-a447b6b360bb.8c7f9b6809be6760e3ef2e19a4c99ce55b21bfff.zh-cn.xlf" "0dab1796-b58b-47e1-aafe-a447b6b360bb.8c7f9b6809be6760e3ef2e19a4c99ce55b21bfff.zh-cn.xlf"
Add-HandbackLink $wsZh "F3" "https://github.com/OpenLocalizationTest/oltest/blob/ab42d9a503830687685b01245f3dbd9aa1fb2bbb/e2e/e94c32b8-8cfd-4625-92ff-3b5067d34b74.md" "e94c32b8-8cfd-4625-92ff-3b5067d34b74.md"
Add-HandbackLink $wsZh "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4415d4333c52fc7d83a3c74ab476dfc18f340275/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/e94c32b8-8cfd-4625-92ff-3b5067d34b74.8af0db742e0519cd28f566026ba04a62020fbfe5.zh-cn.xlf" "e94c32b8-8cfd-4625-92ff-3b5067d34b74.8af0db742e0519cd28f566026ba04a62020fbfe5.zh-cn.xlf"

# --- de-de sheet ---
Add-HandbackLink $wsDe "F2" "https://github.com/OpenLocalizationTest/oltest/blob/ab42d9a503830687685b01245f3dbd9aa1fb2bbb/e2e/0dab1796-b58b-47e1-aafe-a447b6b360bb.md" "0dab1796-b58b-47e1-aafe-a447b6b360bb.md"
Add-HandbackLink $wsDe "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d9d88a72b14bda167a8ba47738040ccfc2c5d23/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/0dab1796-b58b-47e1-aafe-a447b6b360bb.8c7f9b6809be6760e3ef2e19a4c99ce55b21bfff.de-de.xlf" "0dab1796-b58b-47e1-aafe-a447b6b360bb.8c7f9b6809be6760e3ef2e19a4c99ce55b21bfff.de-de.xlf"
Add-HandbackLink $wsDe "F3" "https://github.com/OpenLocalizationTest/oltest/blob/ab42d9a503830687685b01245f3dbd9aa1fb2bbb/e2e/e94c32b8-8cfd-4625-92ff-3b5067d34b74.md" "e94c32b8-8cfd-4625-92ff-3b5067d34b74.md"
Add-HandbackLink $wsDe "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d9d88a72b14bda167a8ba47738040ccfc2c5d23/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/e94c32b8-8cfd-4625-92ff-3b5067d34b74.8af0db742e0519cd28f566026ba04a62020fbfe5.de-de.xlf" "e94c32b8-8cfd-4625-92ff-3b5067d34b74.8af0db742e0519cd28f566026ba04a62020fbfe5.de-de.xlf"
